$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formula in A9: B9+C9 -> B9+C9+3 (value becomes 5)
$ws.Range("A9").Formula = "=B9+C9+3"

# Move the selection from A9 to A10
$ws.Range("A10").Select()
